$d = $word.ActiveDocument

# "Trange instead of T50": the TEMPERATURE column (column 3) of the table
# currently holds T50 values; replace them with the corresponding Trange
# values:
#   0.6429 -> 0.4545
#   0.2143 -> 0.5455
$table = $d.Tables.Item(1)

$map = @{ "0.6429" = "0.4545"; "0.2143" = "0.5455" }

for ($r = 2; $r -le $table.Rows.Count; $r++) {
    $cell = $table.Cell($r, 3)
    $range = $cell.Range
    $range.End = $range.End - 1
    $old = $range.Text
    if ($map.ContainsKey($old)) {
        $range.Text = $map[$old]
    }
}
